# Rotate the "Recorded By" (column G) multi-value lists.
# For every cell in column G that contains a comma-separated list of
# recorder names/emails, move the last entry in the list to the front
# (a right-rotation by one), unless the last entry is exactly "System"
# (those lists are already in their canonical order and are left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ", "
    $n = $parts.Count
    if ($n -lt 2) { continue }

    $last = $parts[$n - 1]
    if ($last -eq "System") { continue }

    $rest = $parts[0..($n - 2)]
    $newParts = @($last) + $rest
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
